# Update existing rows 4 and 5 (dates shifted later by one month) and
# append a new row 6 with a new activity entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: shift start/end dates forward
$ws.Range("B4").Value = (Get-Date -Year 2019 -Month 9 -Day 23 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C4").Value = (Get-Date -Year 2019 -Month 9 -Day 24 -Hour 0 -Minute 0 -Second 0).Date

# Row 5: shift start/end dates forward
$ws.Range("B5").Value = (Get-Date -Year 2019 -Month 9 -Day 24 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C5").Value = (Get-Date -Year 2019 -Month 9 -Day 25 -Hour 0 -Minute 0 -Second 0).Date

# Row 6: new activity row
$ws.Range("A6").Value = "Research about VAT, Tariff, Duties and Excise Tax"
$ws.Range("B6").Value = (Get-Date -Year 2019 -Month 10 -Day 2 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C6").Value = (Get-Date -Year 2019 -Month 10 -Day 3 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("D6").Value = 0.79166666666666663
$ws.Range("E6").Value = 0.91666666666666663

$ws.Range("B6").NumberFormat = $ws.Range("B5").NumberFormat
$ws.Range("C6").NumberFormat = $ws.Range("C5").NumberFormat
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("E6").NumberFormat = $ws.Range("E5").NumberFormat

# Update the selection to match the final state recorded in the workbook
$ws.Range("D10").Select()
